# UC009 - Prestar Contas workbook update: v1.0.3 -> v1.1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) TC1 step 2 expected results: mention "dados do beneficiário"
$ws.Range("D11").Value = "SYSTEM Exibe os detalhes referentes à solicitação selecionada, dados do beneficiário, bem como identificando e apresentando os tipos de documentos/comprovantes a serem informados/consultados pelo usuário; e Exibe o histórico da tramitação da prestação de contas."

# 2) TC2 step 3 action: "Visualiza" -> "Verifica"
$ws.Range("B22").Value = "Chefe Verifica o histório da tramitação da prestação de contas."

# 3) Swap TC3/TC4 step-3 content:
#    TC3 (row 32) now gets the "detalhar a solicitação de diária" content
#    TC4 (row 41) now gets the "visualizar comprovante" content
$ws.Range("B32").Value = "Chefe Clica para detalhar a solicitação de diária."
$ws.Range("D32").Value = "SYSTEM Apresenta a tela de Detalhar Diárias"
$ws.Range("B41").Value = "Chefe Clica em visualizar comprovante."
$ws.Range("D41").Value = "SYSTEM Exibe modal com o comprovante."

# 4) Swap TC7/TC8 step-2 expected-results content:
#    TC7 (row 67) now gets the "não está em nenhum desses dois estados" text
#    TC8 (row 75) now gets the "ainda não pode ter sua prestação de contas realizada" text
$ws.Range("D67").Value = "SYSTEM Identifica que a prestação de contas indicada pelo usuário não está em nenhum desses dois estados: a) NÃO REALIZADA e b) DEVOLVIDA; Permite não permite um novo envio ou alterações na prestação (exclusão de documentos)."
$ws.Range("D75").Value = "SYSTEM Identifica que a solicitação indicada pelo usuário ainda não pode ter sua prestação de contas realizada; Exibe mensagem de erro (MSG212 - Prestação de contas ainda não pode ser realizada) para o usuário, impedindo que ele preste contas (anexa arquivos e etc)."
